$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2022-12-12 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2022-12-13 Tuesday", 2) | Out-Null

# Update every cell in the table with the new arithmetic problems
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "14-7="
$t.Cell(1,2).Range.Text = "44-15="
$t.Cell(1,3).Range.Text = "3+45="
$t.Cell(1,4).Range.Text = "11-3="
$t.Cell(1,5).Range.Text = "87+8="

$t.Cell(2,1).Range.Text = "8+0="
$t.Cell(2,2).Range.Text = "12-10="
$t.Cell(2,3).Range.Text = "9+85="
$t.Cell(2,4).Range.Text = "15+19="
$t.Cell(2,5).Range.Text = "6+62="

$t.Cell(3,1).Range.Text = "69-16="
$t.Cell(3,2).Range.Text = "97-1="
$t.Cell(3,3).Range.Text = "32+65="
$t.Cell(3,4).Range.Text = "71-20="
$t.Cell(3,5).Range.Text = "35-13="

$t.Cell(4,1).Range.Text = "86-4="
$t.Cell(4,2).Range.Text = "16+79="
$t.Cell(4,3).Range.Text = "14+67="
$t.Cell(4,4).Range.Text = "37+29="
$t.Cell(4,5).Range.Text = "30-19="

$t.Cell(5,1).Range.Text = "4-0="
$t.Cell(5,2).Range.Text = "61-2="
$t.Cell(5,3).Range.Text = "93-48="
$t.Cell(5,4).Range.Text = "48+18="
$t.Cell(5,5).Range.Text = "5+7="

$t.Cell(6,1).Range.Text = "16-8="
$t.Cell(6,2).Range.Text = "77-4="
$t.Cell(6,3).Range.Text = "56-18="
$t.Cell(6,4).Range.Text = "85-1="
$t.Cell(6,5).Range.Text = "4+63="

$t.Cell(7,1).Range.Text = "29-23="
$t.Cell(7,2).Range.Text = "19+54="
$t.Cell(7,3).Range.Text = "38-0="
$t.Cell(7,4).Range.Text = "24+64="
$t.Cell(7,5).Range.Text = "3+88="

$t.Cell(8,1).Range.Text = "13+75="
$t.Cell(8,2).Range.Text = "0+45="
$t.Cell(8,3).Range.Text = "79-3="
$t.Cell(8,4).Range.Text = "90-38="
$t.Cell(8,5).Range.Text = "39+23="

$t.Cell(9,1).Range.Text = "26+7="
$t.Cell(9,2).Range.Text = "28-18="
$t.Cell(9,3).Range.Text = "2+21="
$t.Cell(9,4).Range.Text = "30-4="
$t.Cell(9,5).Range.Text = "50-23="

$t.Cell(10,1).Range.Text = "70-49="
$t.Cell(10,2).Range.Text = "27+37="
$t.Cell(10,3).Range.Text = "44-3="
$t.Cell(10,4).Range.Text = "26+27="
$t.Cell(10,5).Range.Text = "56-1="

$t.Cell(11,1).Range.Text = "97-20="
$t.Cell(11,2).Range.Text = "98-43="
$t.Cell(11,3).Range.Text = "72+11="
$t.Cell(11,4).Range.Text = "97-11="
$t.Cell(11,5).Range.Text = "19+54="

$t.Cell(12,1).Range.Text = "15+54="
$t.Cell(12,2).Range.Text = "18+38="
$t.Cell(12,3).Range.Text = "49-14="
$t.Cell(12,4).Range.Text = "0+51="
$t.Cell(12,5).Range.Text = "98-59="

$t.Cell(13,1).Range.Text = "49+40="
$t.Cell(13,2).Range.Text = "3+76="
$t.Cell(13,3).Range.Text = "54+32="
$t.Cell(13,4).Range.Text = "43-21="
$t.Cell(13,5).Range.Text = "91-59="

$t.Cell(14,1).Range.Text = "32+63="
$t.Cell(14,2).Range.Text = "36-33="
$t.Cell(14,3).Range.Text = "94-39="
$t.Cell(14,4).Range.Text = "28+5="
$t.Cell(14,5).Range.Text = "25+61="

$t.Cell(15,1).Range.Text = "69-2="
$t.Cell(15,2).Range.Text = "35+2="
$t.Cell(15,3).Range.Text = "41+28="
$t.Cell(15,4).Range.Text = "55-11="
$t.Cell(15,5).Range.Text = "51+38="

$t.Cell(16,1).Range.Text = "86-46="
$t.Cell(16,2).Range.Text = "83-19="
$t.Cell(16,3).Range.Text = "77-72="
$t.Cell(16,4).Range.Text = "43+17="
$t.Cell(16,5).Range.Text = "97-38="

$t.Cell(17,1).Range.Text = "30+27="
$t.Cell(17,2).Range.Text = "49-19="
$t.Cell(17,3).Range.Text = "76+8="
$t.Cell(17,4).Range.Text = "61+30="
$t.Cell(17,5).Range.Text = "84-81="

$t.Cell(18,1).Range.Text = "90-2="
$t.Cell(18,2).Range.Text = "42+38="
$t.Cell(18,3).Range.Text = "92-78="
$t.Cell(18,4).Range.Text = "82-11="
$t.Cell(18,5).Range.Text = "77-75="

$t.Cell(19,1).Range.Text = "82-20="
$t.Cell(19,2).Range.Text = "76-30="
$t.Cell(19,3).Range.Text = "95-80="
$t.Cell(19,4).Range.Text = "27-25="
$t.Cell(19,5).Range.Text = "97-20="

$t.Cell(20,1).Range.Text = "64+35="
$t.Cell(20,2).Range.Text = "68+11="
$t.Cell(20,3).Range.Text = "52-24="
$t.Cell(20,4).Range.Text = "65+22="
$t.Cell(20,5).Range.Text = "83-55="
